$d = $word.ActiveDocument
$d.Content.Find.Execute("Team ID", $true, $false, $false, $false, $false, $true, 1, $false, "Team", 2)
